# S14/G01: Advanced order types and stop-loss controls
# Appends new sprint-task rows (112-119) to the tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 112; A = "S14"; B = "G01"; C = "Advanced order types and stop-loss controls"; D = "S14_G01_TB001"; E = "Extend Order model and APIs to support Zerodha order types MARKET / LIMIT / SL / SL-M plus trigger price and trigger percent fields."; F = "trigger_percent will be interpreted relative to the Zerodha last traded price (LTP), not the current limit price, and saved alongside trigger_price."; G = "pending"; I = "Requires careful validation so SL/SL-M cannot be placed on the wrong side of the market for BUY/SELL." },
    @{ Row = 113; A = "S14"; B = "G01"; C = "Advanced order types and stop-loss controls"; D = "S14_G01_TB002"; E = "Update execute_order to route SL and SL-M correctly to KiteConnect (trigger_price mandatory, price optional for SL-M) and add guardrails for valid stop-loss placement."; F = "Guardrails include checking trigger_price vs LTP and direction (BUY stops below market, SELL stops above) and rejecting obviously invalid combinations with clear error messages."; G = "pending"; I = "Depends on S14_G01_TB001; behaviour should be documented clearly in PRD and user docs before enabling AUTO strategies to use SL/SL-M." },
    @{ Row = 114; A = "S14"; B = "G02"; C = "Funds and margin preview for edited orders"; D = "S14_G02_TB001"; E = "Wrap Zerodha margins and order_margins APIs in the backend and expose endpoints to fetch available funds and a margin/charges preview for a hypothetical order."; F = "Focus first on the equity segment for Zerodha; later sprints can extend to derivatives or other brokers."; G = "pending"; I = "These endpoints will be used by the queue edit dialog to show Required vs Available amounts including mandatory taxes and charges." },
    @{ Row = 115; A = "S14"; B = "G02"; C = "Funds and margin preview for edited orders"; D = "S14_G02_TF002"; E = "Enhance the Waiting Queue edit dialog to display Required funds (incl. charges) and Available funds, updating dynamically as qty/price/type change."; F = "Uses the new Zerodha preview APIs to recompute required margin after each significant edit, with warnings when funds are insufficient."; G = "pending"; I = "Keep the UI lightweight and responsive; avoid blocking edits if the preview API is temporarily unavailable." },
    @{ Row = 116; A = "S14"; B = "G03"; C = "Queue edit UX polish and stop-loss helpers"; D = "S14_G03_TF001"; E = "Refactor the edit queue order dialog into clear sections (Quantity & price, Stop-loss, Product & preferences, Funds) and add fields for trigger price and trigger percent."; F = "Trigger percent will be entered as a positive or negative percentage relative to LTP and used to derive trigger_price with inline preview."; G = "pending"; I = "Ensure validation and copy make it hard to misconfigure SL/SL-M orders, especially when prices move between alert and execution." },
    @{ Row = 117; A = "S15"; B = "G01"; C = "Zerodha GTT order support"; D = "S15_G01_TB001"; E = "Design how SigmaTrader will map queue orders and preferences into Zerodha GTT single-leg orders (trigger values, last_price source, and order payload)."; F = "GTT design should clarify when to use GTT vs regular orders, how to represent GTT status in SigmaTrader, and how TradingView alerts can request GTT creation."; G = "pending"; I = "This design underpins safe GTT usage for CNC swing trades and must consider off-market placement and modification flows." },
    @{ Row = 118; A = "S15"; B = "G01"; C = "Zerodha GTT order support"; D = "S15_G01_TB002"; E = "Implement backend support for placing, listing, and cancelling Zerodha GTTs using KiteConnect place_gtt / get_gtts / delete_gtt, wired to per-user broker connections."; F = "Initial focus on single-leg GTTs for equity; OCO/advanced patterns can be added later."; G = "pending"; I = "Requires careful error handling and alignment between SigmaTrader order records and Zerodha GTT IDs." },
    @{ Row = 119; A = "S15"; B = "G01"; C = "Zerodha GTT order support"; D = "S15_G01_TF003"; E = "Extend the queue edit and manual order flows to allow creating GTT orders (instead of or in addition to regular orders) when the user selects a GTT option."; F = 'The existing "Convert to GTT" checkbox will be repurposed into a concrete GTT mode that creates or updates real GTTs at Zerodha rather than acting as a passive preference.'; G = "pending"; I = "UI should clearly distinguish between regular orders and GTTs and indicate when an order has an associated active GTT at the broker." }
)

# Columns to populate per row - column H is intentionally left untouched
# (blank / not present), matching the source data for these new rows.
$cols = @(1, 2, 3, 4, 5, 6, 7, 9)
$colKeys = @{ 1 = "A"; 2 = "B"; 3 = "C"; 4 = "D"; 5 = "E"; 6 = "F"; 7 = "G"; 9 = "I" }

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $cell = $ws.Cells.Item($r.Row, $col)
        $cell.Value = $r[$colKeys[$col]]
        # New rows use the workbook default (unwrapped) style rather than the
        # sheet's usual wrap-text column style, matching the source data.
        $cell.Style = "Normal"
    }
}
